{"js": "// Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single\n// practice-sheet table. Each cell holds exactly one run of text such as\n// \"37\u00f78=4, 5\"; the table has 20 rows (4-row stride: data, blank, blank,\n// blank, ...) and 5 columns, so data lives in rows 0, 4, 8, 12, 16.\n//\n// We replace the text cell-by-cell (rather than with one document-wide\n// search/replace) because a couple of the new values happen to collide\n// with old values elsewhere in the table (e.g. \"56\u00f76=9, 2\" is both an\n// original value and a freshly written one), so a global search done in\n// document order could re-match text our own edit just produced.\n// Scoping each Range.search() call to its own table cell's body makes\n// every lookup unambiguous, and using Range.insertText(\"...\", \"Replace\")\n// on the located hit (instead of clearing/re-inserting the whole cell)\n// keeps the existing run/paragraph formatting (fonts, size, alignment)\n// untouched \u2014 only the <w:t> text itself changes, matching the diff.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"items\");\nawait context.sync();\n\n// row, column, expected current text, replacement text\nconst edits = [\n  [0, 0, \"37\u00f78=4, 5\", \"56\u00f76=9, 2\"],\n  [0, 1, \"31\u00f72=15, 1\", \"82\u00f79=9, 1\"],\n  [0, 2, \"86\u00f75=17, 1\", \"99\u00f75=19, 4\"],\n  [0, 3, \"29\u00f79=3, 2\", \"53\u00f76=8, 5\"],\n  [0, 4, \"22\u00f74=5, 2\", \"43\u00f78=5, 3\"],\n\n  [4, 0, \"85\u00f79=9, 4\", \"63\u00f79=7, 0\"],\n  [4, 1, \"58\u00f79=6, 4\", \"70\u00f74=17, 2\"],\n  [4, 2, \"14\u00f79=1, 5\", \"25\u00f79=2, 7\"],\n  [4, 3, \"74\u00f77=10, 4\", \"56\u00f78=7, 0\"],\n  [4, 4, \"96\u00f79=10, 6\", \"85\u00f76=14, 1\"],\n\n  [8, 0, \"26\u00f73=8, 2\", \"57\u00f77=8, 1\"],\n  [8, 1, \"56\u00f76=9, 2\", \"88\u00f74=22, 0\"],\n  [8, 2, \"64\u00f79=7, 1\", \"78\u00f75=15, 3\"],\n  [8, 3, \"21\u00f75=4, 1\", \"90\u00f73=30, 0\"],\n  [8, 4, \"59\u00f75=11, 4\", \"73\u00f77=10, 3\"],\n\n  [12, 0, \"78\u00f73=26, 0\", \"48\u00f72=24, 0\"],\n  [12, 1, \"63\u00f78=7, 7\", \"87\u00f75=17, 2\"],\n  [12, 2, \"30\u00f76=5, 0\", \"56\u00f72=28, 0\"],\n  [12, 3, \"68\u00f77=9, 5\", \"38\u00f73=12, 2\"],\n  [12, 4, \"80\u00f73=26, 2\", \"90\u00f76=15, 0\"],\n\n  [16, 0, \"29\u00f74=7, 1\", \"98\u00f74=24, 2\"],\n  [16, 1, \"78\u00f77=11, 1\", \"77\u00f79=8, 5\"],\n  [16, 2, \"64\u00f74=16, 0\", \"16\u00f79=1, 7\"],\n  [16, 3, \"27\u00f76=4, 3\", \"84\u00f76=14, 0\"],\n  [16, 4, \"84\u00f77=12, 0\", \"77\u00f79=8, 5\"],\n];\n\nfor (const [row, col, oldText, newText] of edits) {\n  const cell = table.getCell(row, col);\n\n  let hits = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length === 0) {\n    // Defensive fallback (shouldn't trigger against the documented source\n    // file): retry case-insensitively before giving up, so that a minor\n    // text variation still gets a targeted, formatting-preserving\n    // replacement instead of falling back to clearing the whole cell.\n    hits = cell.body.search(oldText, { matchCase: false, matchWholeWord: false });\n    hits.load(\"items\");\n    await context.sync();\n  }\n\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    cell.body.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single\n# practice-sheet table. The table has 20 rows (stride of 4: one data row\n# followed by 3 blank rows) and 5 columns, so the populated rows are\n# 1, 5, 9, 13, 17 (Word COM collections are 1-indexed).\n#\n# NOTE: we deliberately do NOT use Content.Find.Execute(..., Replace:=\n# wdReplaceAll) here. In this host, Find.Execute replaces every matching\n# occurrence in the whole document story, regardless of which Range's\n# .Find object invoked it (it is not confined to the calling Range) -\n# and several of the new values we're writing happen to equal *other*\n# cells' original values elsewhere in the table (e.g. \"56\u00f76=9, 2\" is\n# both an original value and a value we write fresh into a different\n# cell), so a text-search-based replace can clobber the wrong cell once\n# a duplicate exists mid-run. Addressing each cell directly by\n# (row, column) and assigning Range.Text sidesteps ambiguity entirely\n# while still preserving the existing run/paragraph formatting (fonts,\n# size, alignment) of each cell, since only the text inside the cell's\n# range is replaced.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# row, column (1-indexed), expected current text, new text\n$edits = @(\n  @(1, 1, \"37\u00f78=4, 5\", \"56\u00f76=9, 2\"),\n  @(1, 2, \"31\u00f72=15, 1\", \"82\u00f79=9, 1\"),\n  @(1, 3, \"86\u00f75=17, 1\", \"99\u00f75=19, 4\"),\n  @(1, 4, \"29\u00f79=3, 2\", \"53\u00f76=8, 5\"),\n  @(1, 5, \"22\u00f74=5, 2\", \"43\u00f78=5, 3\"),\n\n  @(5, 1, \"85\u00f79=9, 4\", \"63\u00f79=7, 0\"),\n  @(5, 2, \"58\u00f79=6, 4\", \"70\u00f74=17, 2\"),\n  @(5, 3, \"14\u00f79=1, 5\", \"25\u00f79=2, 7\"),\n  @(5, 4, \"74\u00f77=10, 4\", \"56\u00f78=7, 0\"),\n  @(5, 5, \"96\u00f79=10, 6\", \"85\u00f76=14, 1\"),\n\n  @(9, 1, \"26\u00f73=8, 2\", \"57\u00f77=8, 1\"),\n  @(9, 2, \"56\u00f76=9, 2\", \"88\u00f74=22, 0\"),\n  @(9, 3, \"64\u00f79=7, 1\", \"78\u00f75=15, 3\"),\n  @(9, 4, \"21\u00f75=4, 1\", \"90\u00f73=30, 0\"),\n  @(9, 5, \"59\u00f75=11, 4\", \"73\u00f77=10, 3\"),\n\n  @(13, 1, \"78\u00f73=26, 0\", \"48\u00f72=24, 0\"),\n  @(13, 2, \"63\u00f78=7, 7\", \"87\u00f75=17, 2\"),\n  @(13, 3, \"30\u00f76=5, 0\", \"56\u00f72=28, 0\"),\n  @(13, 4, \"68\u00f77=9, 5\", \"38\u00f73=12, 2\"),\n  @(13, 5, \"80\u00f73=26, 2\", \"90\u00f76=15, 0\"),\n\n  @(17, 1, \"29\u00f74=7, 1\", \"98\u00f74=24, 2\"),\n  @(17, 2, \"78\u00f77=11, 1\", \"77\u00f79=8, 5\"),\n  @(17, 3, \"64\u00f74=16, 0\", \"16\u00f79=1, 7\"),\n  @(17, 4, \"27\u00f76=4, 3\", \"84\u00f76=14, 0\"),\n  @(17, 5, \"84\u00f77=12, 0\", \"77\u00f79=8, 5\")\n)\n\nforeach ($edit in $edits) {\n  $row = $edit[0]\n  $col = $edit[1]\n  $oldText = $edit[2]\n  $newText = $edit[3]\n\n  $cellRange = $t.Cell($row, $col).Range\n  # Range.Text for a cell includes the trailing cell-mark character(s);\n  # strip them off to compare against the plain expected value.\n  $currentText = $cellRange.Text.TrimEnd([char]7, [char]13)\n  if ($currentText -ne $oldText) {\n    Write-Output (\"WARNING: cell ($row,$col) was \" + $currentText + \", expected \" + $oldText)\n  }\n\n  $cellRange.Text = $newText\n}\n"}
